# Modificacion en psp'p SC_Lectura
#
# Mirrors the upstream change to
# "Documentación/Psp's/Tania/Clase SC_Lectura/MetricasGenerales.xlsx":
#   - the "Report generated at ..." footer text was refreshed to a later run
#   - most of the Plan/Actual/To-Date metric numbers were updated to match
#     that later report run
#   - the saved window size in the workbook view changed (best effort - the
#     host may not persist this cosmetic, non-data setting)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Footer timestamp text (A27) ---------------------------------------
$ws.Range("A27").Value = "Report generated at 9:52 PM on Dec 6, 2018"

# --- Updated metric values ----------------------------------------------
# Row 8 - Size/Hour
$ws.Range("C8").Value  = 8.39
$ws.Range("D8").Value  = 22.7

# Row 10 - Actual Time ([h]:mm)
$ws.Range("D10").Value = 0.41875000000000001

# Row 11 - CPI (Cost-Performance Index)
$ws.Range("D11").Value = 1.19

# Row 14 - Test Defects/KLOC or equivalent
$ws.Range("C14").Value = 44.4
$ws.Range("D14").Value = 17.5

# Row 15 - Total Defects/KLOC or equivalent
$ws.Range("C15").Value = 66.7
$ws.Range("D15").Value = 35.1

# Row 16 - Yield %
$ws.Range("D16").Value = 0.5

# Row 17 - Code Review Rate
$ws.Range("C17").Value = 33.799999999999997
$ws.Range("D17").Value = 130

# Row 18 - % Appraisal COQ
$ws.Range("D18").Value = 0.27400000000000002

# Row 19 - % Failure COQ
$ws.Range("D19").Value = 0.14399999999999999

# Row 20 - COQ A/F Ratio
$ws.Range("D20").Value = 1.9

# Row 21 - PQI
$ws.Range("D21").Value = 0.23

# --- Saved window size (best effort; cosmetic workbook-view setting) ---
$win = $excel.ActiveWindow
$win.Width = 8790
$win.Height = 6750
